# Rename the "adj_r_squared" sheet to "metrics" and replace its single
# Adj.R^2 value with a small metrics table (name/value pairs):
#   metrics   | value
#   Adj.R^2   | 0.663618471219521
#   NRMSE     | 0.0210302853049309
#   SMAPE     | 0.0163438415539451
#   RMSE      | 0.214155864195739

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("adj_r_squared")

$ws.Name = "metrics"

$ws.Range("A1").Value = "metrics"
$ws.Range("B1").Value = "value"

$ws.Range("A2").Value = "Adj.R^2"
$ws.Range("B2").Value = 0.663618471219521

$ws.Range("A3").Value = "NRMSE"
$ws.Range("B3").Value = 0.0210302853049309

$ws.Range("A4").Value = "SMAPE"
$ws.Range("B4").Value = 0.0163438415539451

$ws.Range("A5").Value = "RMSE"
$ws.Range("B5").Value = 0.214155864195739
